# Update calibration estimates across the "AR", "SETAR", "GARCH", "TARCH"
# and "AR_TARCH" sheets with new values computed using a std moving average.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.005476268475053387
$ws.Range("B3").Value = 0.7708829775538072
$ws.Range("B4").Value = 0.09883757256385761
$ws.Range("B5").Value = "[0.9999999999999999, -0.0029305740574956055, -0.06430183407658344, -0.09678776685102415, -0.11332665377180157, 0.060223909508484656, 0.023134437068713504, 0.027527342155996462, 0.021708550479818258, 0.0037143385741970026, -0.037151165649365075, -0.04516230383121428, -0.03225842474746748, 0.028619202486163828, 0.05508668071366697, 0.03165893787633444, -0.03250348169091397, -0.04772375330653062, -0.0384930840249508, 0.009568695150526942]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1792157409623108
$ws.Range("B5").Value = 0.5881546956886264
$ws.Range("B6").Value = 0.05512821606360369
$ws.Range("B7").Value = 0.1612931945760407
$ws.Range("B8").Value = 0.5596449676183822
$ws.Range("B9").Value = 0.05645687639108358
$ws.Range("B10").Value = "[0.9999999999999998, 0.045570465761380005, 0.012460174431676737, -0.017206711107716557, 0.03238529985621976, 0.011684009901597903, -0.002676436387977331, 0.04309004426792053, -0.010217509357506123, -0.009565570123857064, -0.0028184498406111738, -0.022414792035762885, -0.013905355321581497, 0.0038708339911607548, -0.004739660579801487, 0.006849143917652729, -0.03921815225943939, -0.031256743219317415, -0.030111662746606812, 0.00072065700403912]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.00009402145416688678
$ws.Range("B3").Value = 0.002485355801549447
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0.9774384499711841
$ws.Range("B6").Value = "[1.0, 0.0053512265578218865, -0.047965270758115824, -0.09624917996587308, -0.12111843580433812, 0.11133897598606488, 0.028929559458363084, 0.02794793613062658, 0.030869708830733384, 0.0057766687575217576, -0.06406924911461513, -0.045401848281271855, -0.030071506666856523, 0.0364771871233387, 0.032310908958978965, 0.024703091277189164, -0.04091805459745477, -0.059298038338594115, -0.04759210432023713, 0.012685546856744022]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.0005872119705266554
$ws.Range("B3").Value = 0.1108284876514284
$ws.Range("B4").Value = 0.01473822355044238
$ws.Range("B5").Value = 0.00000000008254491197321024
$ws.Range("B6").Value = "[0.9999999999999998, -0.0006420118170699893, -0.04414728998320402, -0.09554911476183414, -0.11878321499340518, 0.11447786515898213, 0.0264906220059799, 0.029747048250579198, 0.029623072948371954, 0.008194965448504015, -0.06251944152304799, -0.043928346907034316, -0.02838457633446304, 0.036959261044889564, 0.03444569340507933, 0.024917525691859312, -0.03859896143881191, -0.05734226218217965, -0.0451869741055887, 0.013668181538868819]"
$ws.Range("B7").Value = -0.01473822355536685

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.005455171079608237
$ws.Range("B3").Value = 0.07055268146308906
$ws.Range("B4").Value = 0.0000002297793494844305
$ws.Range("B5").Value = 0.2860445153941196
$ws.Range("B6").Value = "[0.9999999999999998, -0.002851912227737137, -0.06427066988832664, -0.09673299731561498, -0.1133420878520739, 0.059996615575479764, 0.023155539201167287, 0.027443463001512357, 0.021780452573349072, 0.003759274871357114, -0.03717623372531297, -0.04509679349928532, -0.03244773821903248, 0.02867704038661088, 0.05501713893550652, 0.03171998396605875, -0.032465726486331455, -0.04779459930836689, -0.038439169918806304, 0.00941241838073797]"
$ws.Range("B7").Value = -0.0000002207543423280119
$ws.Range("B9").Value = 0.7708307639700828
